$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 0.13
$ws.Range("I2").Value = 0.42
$ws.Range("J2").Value = 0.114
$ws.Range("K2").Value = 0.39
$ws.Range("H3").Value = 0.871
$ws.Range("J3").Value = 0.771
$ws.Range("K3").Value = 0.9399999999999999
$ws.Range("H4").Value = 0.302
$ws.Range("I4").Value = 0.6
$ws.Range("J4").Value = 0.222
$ws.Range("K4").Value = 0.54
$ws.Range("G5").Value = 0.38
$ws.Range("H5").Value = 0.529
$ws.Range("I5").Value = 0.82
$ws.Range("J5").Value = 0.317
$ws.Range("K5").Value = 0.6899999999999999
$ws.Range("E6").Value = 0.31
$ws.Range("G6").Value = 0.47
$ws.Range("H6").Value = 0.532
$ws.Range("I6").Value = 0.58
$ws.Range("J6").Value = 0.3
$ws.Range("K6").Value = 0.6
$ws.Range("E7").Value = 0.48
$ws.Range("G7").Value = 0.48
$ws.Range("H7").Value = 0.783
$ws.Range("I7").Value = 0.93
$ws.Range("J7").Value = 0.516
$ws.Range("K7").Value = 0.83
$ws.Range("E8").Value = 0.2
$ws.Range("G8").Value = 0.36
$ws.Range("H8").Value = 0.454
$ws.Range("I8").Value = 0.62
$ws.Range("J8").Value = 0.388
$ws.Range("K8").Value = 0.63
$ws.Range("G9").Value = 0.22
$ws.Range("H9").Value = 0.154
$ws.Range("I9").Value = 0.47
$ws.Range("J9").Value = 0.094
$ws.Range("K9").Value = 0.34
$ws.Range("G10").Value = 0.2
$ws.Range("H10").Value = 0.133
$ws.Range("I10").Value = 0.42
$ws.Range("J10").Value = 0.07199999999999999
$ws.Range("K10").Value = 0.28
$ws.Range("E12").Value = 0.54
$ws.Range("H12").Value = 0.542
$ws.Range("I12").Value = 0.66
$ws.Range("J12").Value = 0.466
$ws.Range("K12").Value = 0.67
$ws.Range("E13").Value = 0.7
$ws.Range("H13").Value = 0.626
$ws.Range("I13").Value = 0.88
$ws.Range("J13").Value = 0.492
$ws.Range("K13").Value = 0.82
$ws.Range("G14").Value = 0.4
$ws.Range("H14").Value = 0.256
$ws.Range("I14").Value = 0.58
$ws.Range("J14").Value = 0.154
$ws.Range("K14").Value = 0.46
$ws.Range("G16").Value = 0.09
$ws.Range("E18").Value = 0.26
$ws.Range("G18").Value = 0.32
$ws.Range("H18").Value = 0.73
$ws.Range("I18").Value = 0.63
$ws.Range("J18").Value = 0.512
$ws.Range("K18").Value = 0.7
$ws.Range("E20").Value = 0.52
$ws.Range("G20").Value = 0.5600000000000001
$ws.Range("H20").Value = 0.637
$ws.Range("I20").Value = 0.52
$ws.Range("J20").Value = 0.508
$ws.Range("K20").Value = 0.61
$ws.Range("E21").Value = 0.84
$ws.Range("G21").Value = 0.73
$ws.Range("H21").Value = 0.641
$ws.Range("I21").Value = 0.86
$ws.Range("J21").Value = 0.58
$ws.Range("K21").Value = 0.84
$ws.Range("E22").Value = 0.05
$ws.Range("G22").Value = 0.05
$ws.Range("H24").Value = 0.337
$ws.Range("I24").Value = 0.63
$ws.Range("J24").Value = 0.167
$ws.Range("K24").Value = 0.48
$ws.Range("H25").Value = 0.112
$ws.Range("I25").Value = 0.39
$ws.Range("J25").Value = 0.08
$ws.Range("K25").Value = 0.3
$ws.Range("G26").Value = 0.3
$ws.Range("H26").Value = 0.331
$ws.Range("I26").Value = 0.61
$ws.Range("J26").Value = 0.182
$ws.Range("K26").Value = 0.5
$ws.Range("H27").Value = 0.082
$ws.Range("I27").Value = 0.31
$ws.Range("J27").Value = 0.034
$ws.Range("K27").Value = 0.15
$ws.Range("E28").Value = 0.36
$ws.Range("G28").Value = 0.42
$ws.Range("H28").Value = 0.8149999999999999
$ws.Range("I28").Value = 0.38
$ws.Range("J28").Value = 0.498
$ws.Range("K28").Value = 0.65
$ws.Range("E29").Value = 0.35
$ws.Range("G29").Value = 0.3
$ws.Range("H29").Value = 0.852
$ws.Range("I29").Value = 0.92
$ws.Range("J29").Value = 0.649
$ws.Range("K29").Value = 0.87
$ws.Range("H32").Value = 0.248
$ws.Range("I32").Value = 0.62
$ws.Range("J32").Value = 0.257
$ws.Range("K32").Value = 0.63
$ws.Range("H33").Value = 0.129
$ws.Range("I33").Value = 0.42
$ws.Range("J33").Value = 0.095
$ws.Range("K33").Value = 0.34
$ws.Range("H34").Value = 0.242
$ws.Range("I34").Value = 0.62
$ws.Range("J34").Value = 0.212
$ws.Range("K34").Value = 0.57
$ws.Range("H35").Value = 0.378
$ws.Range("I35").Value = 0.75
$ws.Range("J35").Value = 0.257
$ws.Range("K35").Value = 0.63
$ws.Range("H36").Value = 0.594
$ws.Range("I36").Value = 0.88
$ws.Range("J36").Value = 0.393
$ws.Range("K36").Value = 0.76
$ws.Range("H37").Value = 0.095
$ws.Range("I37").Value = 0.35
$ws.Range("J37").Value = 0.059
$ws.Range("K37").Value = 0.24
$ws.Range("G42").Value = 0.65
$ws.Range("H42").Value = 0.167
$ws.Range("I42").Value = 0.5
$ws.Range("J42").Value = 0.115
$ws.Range("K42").Value = 0.39
$ws.Range("H43").Value = 0.308
$ws.Range("I43").Value = 0.6899999999999999
$ws.Range("J43").Value = 0.163
$ws.Range("H44").Value = 0.022
$ws.Range("I44").Value = 0.1
$ws.Range("J44").Value = 0.011
$ws.Range("K44").Value = 0.05
$ws.Range("E45").Value = 0.35
$ws.Range("H45").Value = 0.586
$ws.Range("I45").Value = 0.87
$ws.Range("J45").Value = 0.445
$ws.Range("K45").Value = 0.8
